# Remove the "online" cost data columns (2014 / 2013, columns H:I) and their
# associated per-year summary rows (which lived further right and get
# shifted left automatically when the columns are deleted).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete columns H and I entirely; everything to the right (the Year / Mean /
# Median / High / Low / Count summary table) shifts left by two columns.
$ws.Range("H1:I1").EntireColumn.Delete()

# The last two rows of the summary table (2014 / 2013), now at I9:N10, are
# also removed since their source data (old columns H/I) is gone.
$ws.Range("I9:N10").ClearContents()

# Update the remembered selection to match the authored change.
$ws.Range("Q3").Select()
